$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values for two rows in both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 10167
    $ws.Range("F5").Value = 627
}
